$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$excel.ActiveWindow.ScrollRow = 70
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D95").Select()
